$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RVL")

# Row 7: DDT - DataSources
$ws.Range("A7").Value = ""
$ws.Range("B7").Value = "Action"
$ws.Range("C7").Value = "RVL"
$ws.Range("D7").Value = "DoPlayScript"
$ws.Range("E7").Value = "scriptPath"
$ws.Range("F7").Value = "string"
$ws.Range("G7").Value = "%WORKDIR%\DataSources\Main.rvl.xlsx"

# Row 8: DDT - DataOrigin
$ws.Range("B8").Value = "Action"
$ws.Range("C8").Value = "RVL"
$ws.Range("D8").Value = "DoPlayScript"
$ws.Range("E8").Value = "scriptPath"
$ws.Range("F8").Value = "string"
$ws.Range("G8").Value = "%WORKDIR%\DataOrigin\Main.rvl.xlsx"

# Row 9: DDT - DataOutput
$ws.Range("B9").Value = "Action"
$ws.Range("C9").Value = "RVL"
$ws.Range("D9").Value = "DoPlayScript"
$ws.Range("E9").Value = "scriptPath"
$ws.Range("F9").Value = "string"
$ws.Range("G9").Value = "%WORKDIR%\DataOutput\Main.rvl.xlsx"
